# Update crypto price/volume data per Sat Nov  9 20:59:39 UTC 2024 GitHub Actions run
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "76.291.16"
$ws.Range("E2").Value = "  -0.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.083.73"
$ws.Range("E3").Value = "  +5.12%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "198.41"
$ws.Range("E5").Value = "  +0.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "619.66"
$ws.Range("E6").Value = "  +4.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.551"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("E9").Value = "  +6.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.080.34"
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.442"
$ws.Range("E11").Value = "  +0.48%  "
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("E13").Value = "  +7.81%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.647.66"
$ws.Range("E14").Value = "  +5.13%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "29.31"
$ws.Range("E15").Value = "  +3.37%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.192.62"
$ws.Range("E16").Value = "  -0.39%  "
$ws.Range("E17").Value = "  +3.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.073.32"
$ws.Range("E18").Value = "  +4.95%  "
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.99"
$ws.Range("E20").Value = "  +3.72%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.61"
$ws.Range("E21").Value = "  +15.76%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "381.88"
$ws.Range("E22").Value = "  +2.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.46"
$ws.Range("E23").Value = "  +4.02%  "
$ws.Range("B24").Value = "LEO"
$ws.Range("C24").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.44"
$ws.Range("E24").Value = "  +0.53%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.236.39"
$ws.Range("E25").Value = "  +4.94%  "
$ws.Range("B26").Value = "NEARProtocol"
$ws.Range("C26").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "4.53"
$ws.Range("E26").Value = "  +6.64%  "
$ws.Range("B27").Value = "Litecoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "72.37"
$ws.Range("E27").Value = "  +0.37%  "
$ws.Range("B28").Value = "Dai"
$ws.Range("C28").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.09"
$ws.Range("E29").Value = "  +5.35%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0000108"
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.20%  "
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.29"
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.42"
$ws.Range("E33").Value = "  +3.98%  "
$ws.Range("B34").Value = "Bittensor"
$ws.Range("C34").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "504.08"
$ws.Range("E34").Value = "  +0.33%  "
$ws.Range("B35").Value = "PancakeSwap"
$ws.Range("C35").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.92"
$ws.Range("E35").Value = "  +5.56%  "
$ws.Range("B36").Value = "FirstDigitalUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  -0.02%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.125"
$ws.Range("E37").Value = "  +14.18%  "
$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.83"
$ws.Range("E38").Value = "  +3.54%  "
$ws.Range("B39").Value = "Monero"
$ws.Range("C39").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "162.06"
$ws.Range("E39").Value = "  -2.14%  "
$ws.Range("B40").Value = "Aave"
$ws.Range("C40").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "195.54"
$ws.Range("E40").Value = "  +8.96%  "
$ws.Range("B41").Value = "WhiteBITCoin"
$ws.Range("C41").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "20.06"
$ws.Range("E41").Value = "  +0.54%  "
$ws.Range("B42").Value = "PolygonEcosystemToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.379"
$ws.Range("E42").Value = "  -3.27%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.103"
$ws.Range("E43").Value = "  -7.73%  "
$ws.Range("B44").Value = "USDe"
$ws.Range("C44").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.808"
$ws.Range("E45").Value = "  +23.72%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "5.16"
$ws.Range("E46").Value = "  +5.12%  "
$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.25"
$ws.Range("E47").Value = "  +6.20%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.65"
$ws.Range("E48").Value = "  +0.69%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.44"
$ws.Range("E49").Value = "  +5.97%  "
$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "40.43"
$ws.Range("E50").Value = "  +0.78%  "
$ws.Range("B51").Value = "ARBITRUM"
$ws.Range("C51").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.599"
$ws.Range("E51").Value = "  +1.86%  "

Write-Host "Applied cryptos update"
